$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "Procedures" (sheet1)
# ---------------------------------------------------------------------
$wsProc = $wb.Worksheets.Item("Procedures")

# Row 2: shorten the description text and add a new Algorithm entry
$wsProc.Range("D2").Value = "Setups the variables for the engage run"
$wsProc.Range("E2").Value = "While Check_Value < 5`n       Increment Check_Value + 1`n       If Check_Value = 2`n             Set Ex_Range_Distance_1 `n             Ex_Range_Distance_2 to max`n       If Check_Value = 3`n             Set Ex_Range_Distance_1 `n             Ex_Range_Distance_2 to min`n"
$wsProc.Range("E2").WrapText = $true

# Row 3: replace the old (mostly empty) row with the new Failure procedure
# (Column order matters for shared-string allocation order, matching the
# authoring tool: Failure, Call_Failure, None, Input..., Test_Id...)
$wsProc.Range("A3").Value = "Failure"
$wsProc.Range("C3").Value = "Call_Failure"
$wsProc.Range("B3").Value = "None"
$wsProc.Range("D3").Value = "Input a test Id to see if a failure occurred"
$wsProc.Range("E3").Value = "Test_Id Local = Test_Id`nIf Failure_Data_Structure(Tests_Id_Local) = True"
$wsProc.Range("D3").WrapText = $true
$wsProc.Range("E3").WrapText = $true

# Column + row sizing
$wsProc.Columns.Item(5).ColumnWidth = 37.166666666666664
$wsProc.Rows.Item(2).RowHeight = 165
$wsProc.Rows.Item(3).RowHeight = 60

# ---------------------------------------------------------------------
# Sheet "Constants" (sheet5)
# ---------------------------------------------------------------------
$wsConst = $wb.Worksheets.Item("Constants")

$wsConst.Range("A2").Value = "Simulation"
$wsConst.Range("B2").Value = "Hex_Value_Lock"
$wsConst.Range("C2").Value = "Parent_Hex"
$wsConst.Range("D2").Value = $false
$wsConst.Range("E2").Value = 2

$wsConst.Columns.Item(2).ColumnWidth = 14.666666666666666
$wsConst.Columns.Item(3).ColumnWidth = 10.5

# ---------------------------------------------------------------------
# Selection / active-cell bookkeeping (done last, after all data edits)
# ---------------------------------------------------------------------

# "Range Types" (sheet2) - selection only
$wsRange = $wb.Worksheets.Item("Range Types")
$wsRange.Range("D2").Select()

# "Constants" (sheet5) - selection only
$wsConst.Range("E7").Select()

# "Procedures" (sheet1) - selection, then re-activate so it stays the
# tab that is selected when the workbook is saved
$wsProc.Range("D4").Select()
$wsProc.Activate()
